# Update "想去人数" (F column) figures for several exhibitions across sheets,
# matching the regenerated gh-pages data output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 198
$ws1.Range("F6").Value = 335
$ws1.Range("F7").Value = 223
$ws1.Range("F8").Value = 2173
$ws1.Range("F10").Value = 5343
$ws1.Range("F11").Value = 119

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 51

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 51
$ws4.Range("F6").Value = 198
$ws4.Range("F7").Value = 335
$ws4.Range("F8").Value = 223
$ws4.Range("F11").Value = 2173
$ws4.Range("F13").Value = 5343
$ws4.Range("F14").Value = 119
